$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 56499.5
$ws.Range("J3").Value = 56499.5
$ws.Range("L3").Value = 56499.5
$ws.Range("N3").Value = -56727.5
$ws.Range("H53").Value = 168.03847
$ws.Range("I53").Value = 151.3077
$ws.Range("J53").Value = 184.76923
$ws.Range("K53").Value = 151.3077
$ws.Range("L53").Value = 184.76923
$ws.Range("M53").Value = 485.6923
$ws.Range("N53").Value = -1458.76923
$ws.Range("H55").Value = 429.08334
$ws.Range("I55").Value = 196.125
$ws.Range("J55").Value = 895
$ws.Range("K55").Value = 196.125
$ws.Range("L55").Value = 895
$ws.Range("M55").Value = 17.875
$ws.Range("N55").Value = -1323
$ws.Range("H69").Value = 5671.6665
$ws.Range("J69").Value = 3507.5
$ws.Range("L69").Value = 10522.5
$ws.Range("N69").Value = -12270.5
$ws.Range("H72").Value = 5671.6665
$ws.Range("J72").Value = 3507.5
$ws.Range("L72").Value = 31567.5
$ws.Range("N72").Value = -40303.5
$ws.Range("H74").Value = 4578.5713
$ws.Range("I74").Value = 4537.5
$ws.Range("J74").Value = 4633.3335
$ws.Range("K74").Value = 4537.5
$ws.Range("L74").Value = 4633.3335
$ws.Range("M74").Value = -3601.5
$ws.Range("N74").Value = -6505.3335
$ws.Range("H76").Value = 4245.2383
$ws.Range("I76").Value = 3886.6667
$ws.Range("J76").Value = 5141.6665
$ws.Range("K76").Value = 3886.6667
$ws.Range("L76").Value = 5141.6665
$ws.Range("M76").Value = -3571.6667
$ws.Range("N76").Value = -5771.6665
$ws.Range("H77").Value = 4578.5713
$ws.Range("I77").Value = 4537.5
$ws.Range("J77").Value = 4633.3335
$ws.Range("K77").Value = 22687.5
$ws.Range("L77").Value = 23166.6675
$ws.Range("M77").Value = -18007.5
$ws.Range("N77").Value = -32526.6675
$ws.Range("H79").Value = 4245.2383
$ws.Range("I79").Value = 3886.6667
$ws.Range("J79").Value = 5141.6665
$ws.Range("K79").Value = 3886.6667
$ws.Range("L79").Value = 5141.6665
$ws.Range("M79").Value = -2794.6667
$ws.Range("N79").Value = -7325.6665
$ws.Range("H80").Value = 1182.1428
$ws.Range("I80").Value = 1888
$ws.Range("K80").Value = 5664
$ws.Range("M80").Value = -4666
$ws.Range("H83").Value = 1182.1428
$ws.Range("I83").Value = 1888
$ws.Range("K83").Value = 16992
$ws.Range("M83").Value = -12000
$ws.Range("H92").Value = 2537.5
$ws.Range("I92").Value = 2650
$ws.Range("J92").Value = 2200
$ws.Range("K92").Value = 2650
$ws.Range("L92").Value = 2200
$ws.Range("M92").Value = -1402
$ws.Range("N92").Value = -4696
$ws.Range("H96").Value = 2005.2
$ws.Range("I96").Value = 2006.5
$ws.Range("J96").Value = 2000
$ws.Range("K96").Value = 6019.5
$ws.Range("L96").Value = 6000
$ws.Range("M96").Value = -4646.5
$ws.Range("N96").Value = -8746
$ws.Range("H98").Value = 5185.7144
$ws.Range("I98").Value = 5185.7144
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 5185.7144
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -3687.7144
$ws.Range("N98").Value = $null
$ws.Range("H102").Value = 56499.5
$ws.Range("J102").Value = 56499.5
$ws.Range("L102").Value = 56499.5
$ws.Range("N102").Value = -62989.5
$ws.Range("H106").Value = 3800
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 3800
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 3800
$ws.Range("M106").Value = $null
$ws.Range("N106").Value = -5062
$ws.Range("H107").Value = 988.5714
$ws.Range("I107").Value = 988.5714
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 988.5714
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 931.4286
$ws.Range("N107").Value = $null
$ws.Range("H112").Value = 4688.7144
$ws.Range("J112").Value = 1450.6842
$ws.Range("L112").Value = 4352.0526
$ws.Range("N112").Value = -6568.0526
$ws.Range("H113").Value = 2877.9
$ws.Range("J113").Value = 2974.8572
$ws.Range("L113").Value = 2974.8572
$ws.Range("N113").Value = -9482.8572
$ws.Range("H122").Value = 5185.7144
$ws.Range("I122").Value = 5185.7144
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 15557.1432
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -13107.1432
$ws.Range("N122").Value = $null

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 325.125
$ws.Range("I5").Value = 325.125
$ws.Range("K5").Value = 325.125
$ws.Range("M5").Value = -213.125
$ws.Range("H45").Value = 1489.6182
$ws.Range("I45").Value = 1410.119
$ws.Range("K45").Value = 1410.119
$ws.Range("M45").Value = -1033.119
$ws.Range("H108").Value = 41000
$ws.Range("J108").Value = 41000
$ws.Range("L108").Value = 41000
$ws.Range("N108").Value = -48680

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 325.125
$ws.Range("I4").Value = 325.125
$ws.Range("K4").Value = 325.125
$ws.Range("M4").Value = -210.125
$ws.Range("H132").Value = 65000
$ws.Range("J132").Value = 65000
$ws.Range("L132").Value = 65000
$ws.Range("N132").Value = -75120

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 600393.6
$ws.Range("I31").Value = 4536
$ws.Range("J31").Value = 1060055.2
$ws.Range("K31").Value = 4536
$ws.Range("L31").Value = 1060055.2
$ws.Range("M31").Value = -4241
$ws.Range("N31").Value = -1060645.2
$ws.Range("H34").Value = 600393.6
$ws.Range("I34").Value = 4536
$ws.Range("J34").Value = 1060055.2
$ws.Range("K34").Value = 4536
$ws.Range("L34").Value = 1060055.2
$ws.Range("M34").Value = -4334
$ws.Range("N34").Value = -1060459.2
$ws.Range("H59").Value = 22368.732
$ws.Range("I59").Value = 10104
$ws.Range("J59").Value = 23244.785
$ws.Range("K59").Value = 10104
$ws.Range("L59").Value = 23244.785
$ws.Range("M59").Value = -8959
$ws.Range("N59").Value = -25534.785
$ws.Range("H132").Value = 4290.8945
$ws.Range("I132").Value = 4251.5
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 12754.5
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -10224.5
$ws.Range("N132").Value = -20060

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = $null
$ws.Range("N16").Value = $null
$ws.Range("H74").Value = 6750
$ws.Range("I74").Value = 6500
$ws.Range("K74").Value = 19500
$ws.Range("M74").Value = -18439
$ws.Range("H77").Value = 6750
$ws.Range("I77").Value = 6500
$ws.Range("K77").Value = 58500
$ws.Range("M77").Value = -53196
$ws.Range("H132").Value = 2807.4211
$ws.Range("I132").Value = 2800.4
$ws.Range("J132").Value = 2808.4849
$ws.Range("K132").Value = 25203.6
$ws.Range("L132").Value = 25276.3641
$ws.Range("M132").Value = -22673.6
$ws.Range("N132").Value = -30336.3641
$ws.Range("H136").Value = 3597.8235
$ws.Range("I136").Value = 977.2
$ws.Range("J136").Value = 4689.75
$ws.Range("K136").Value = 2931.6
$ws.Range("L136").Value = 14069.25
$ws.Range("M136").Value = 2168.4
$ws.Range("N136").Value = -24269.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 12031.077
$ws.Range("I122").Value = 15144
$ws.Range("J122").Value = 5027
$ws.Range("K122").Value = 45432
$ws.Range("L122").Value = 15081
$ws.Range("M122").Value = -42982
$ws.Range("N122").Value = -19981
$ws.Range("H126").Value = 2835.7896
$ws.Range("I126").Value = 1916.3636
$ws.Range("J126").Value = 4100
$ws.Range("K126").Value = 5749.0908
$ws.Range("L126").Value = 12300
$ws.Range("M126").Value = -3279.0908
$ws.Range("N126").Value = -17240

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4073.6
$ws.Range("I40").Value = 3475.3333
$ws.Range("J40").Value = 6466.6665
$ws.Range("K40").Value = 3475.3333
$ws.Range("L40").Value = 6466.6665
$ws.Range("M40").Value = -3339.3333
$ws.Range("N40").Value = -6738.6665
$ws.Range("H122").Value = 4859.725
$ws.Range("I122").Value = 5164
$ws.Range("J122").Value = 3811.6667
$ws.Range("K122").Value = 15492
$ws.Range("L122").Value = 11435.0001
$ws.Range("M122").Value = -13042
$ws.Range("N122").Value = -16335.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 500
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 500
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 500
$ws.Range("M3").Value = $null
$ws.Range("N3").Value = -728
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").Value = $null
$ws.Range("H100").Value = 2483.8333
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 2483.8333
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 4967.6666
$ws.Range("M100").Value = $null
$ws.Range("N100").Value = -6049.6666
$ws.Range("H122").Value = 1950.5714
$ws.Range("I122").Value = 1800.7273
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 5402.1819
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -2952.1819
$ws.Range("N122").Value = -12400
$ws.Range("H126").Value = 1806.6471
$ws.Range("I126").Value = 1655.2727
$ws.Range("J126").Value = 2084.1667
$ws.Range("K126").Value = 4965.8181
$ws.Range("L126").Value = 6252.500100000001
$ws.Range("M126").Value = -2495.8181
$ws.Range("N126").Value = -11192.5001
$ws.Range("H132").Value = 2206.5278
$ws.Range("I132").Value = 2486.0356
$ws.Range("J132").Value = 1228.25
$ws.Range("K132").Value = 7458.1068
$ws.Range("L132").Value = 3684.75
$ws.Range("M132").Value = -4928.1068
$ws.Range("N132").Value = -8744.75
